# feat/ added new data
#
# Adds a new round entry (Russell @ Kapiti Golf Course) to the bottom of the
# "Summary" log table, then re-sorts the whole data range by Date (column A)
# ascending - mirroring the existing AutoFilter sort state so the table
# stays chronologically ordered after the new row is typed in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- locate the current bottom of the table -------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$newRow = $lastRow + 1

# --- append the new round as a new row at the bottom -----------------------
$ws.Cells.Item($newRow, 1).Value = 46068
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
$ws.Cells.Item($newRow, 2).Value = "Kapiti Golf Course"
$ws.Cells.Item($newRow, 3).Value = "Russell"
$ws.Cells.Item($newRow, 4).Value = "Front-9"
$ws.Cells.Item($newRow, 5).Value = "Ambros"
$ws.Cells.Item($newRow, 6).Value = 47
$ws.Cells.Item($newRow, 7).Value = 34
$ws.Cells.Item($newRow, 8).Formula = "=SUM(F" + $newRow + "-G" + $newRow + ")"
$ws.Cells.Item($newRow, 9).Value = "Not my greatest round haha. Great to get Kat out!"

# --- re-sort the whole data range (excludes header row) by Date asc --------
$sortRange = $ws.Range("A2:I" + $newRow)
$keyRange = $ws.Range("A2:A" + $newRow)
$sortRange.Sort($keyRange, 1)

# --- keep the AutoFilter / selection in sync with the new extent -----------
$ws.Range("A1:I" + $newRow).AutoFilter(1)
$ws.Range("I" + $newRow).Select()
